# Weekly update: insert 3 new daily price rows for "Alcachofa" (Hortaliza,
# Mercado Mayorista Lo Valledor de Santiago) reported on 2023-06-29
# (Excel serial 45106), pushing the existing historical rows down by 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows above row 899; rows 899:967 shift down to 902:970.
$ws.Rows("899:901").Insert()

# Columns: A..R
# A=Mercado ID, B=Mercado, C=Region, D=Fecha, E=Codreg, F=Categoria ID,
# G=Categoria, H=Variedad, I=Calidad, J=Volumen, K=Precio minimo,
# L=Precio maximo, M=Precio promedio ponderado, N=Unidad de comercializacion,
# O=Origen, P=Precio $/Kg, Q=Kg o Unidades, R=Clasificacion
$newRows = @(
    @{ Row=899; H="Argentina(o)"; I="Primera"; J=400;  K=13000; L=14000; M=13425; N="$/caja 50 unidades"; P=268;   Q=50 },
    @{ Row=900; H="Española";     I="Extra";   J=580;  K=17000; L=18000; M=17448; N="$/caja 25 unidades"; P=17448; Q=1  },
    @{ Row=901; H="Española";     I="Primera"; J=1040; K=18000; L=19000; M=18663; N="$/caja 30 unidades"; P=622;   Q=30 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value2 = 6
    $ws.Range("B$row").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
    $ws.Range("C$row").Value2 = "Metropolitana"
    $ws.Range("D$row").Value2 = 45106
    $ws.Range("E$row").Value2 = 13
    $ws.Range("F$row").Value2 = 100112013
    $ws.Range("G$row").Value2 = "Alcachofa"
    $ws.Range("H$row").Value2 = $r.H
    $ws.Range("I$row").Value2 = $r.I
    $ws.Range("J$row").Value2 = $r.J
    $ws.Range("K$row").Value2 = $r.K
    $ws.Range("L$row").Value2 = $r.L
    $ws.Range("M$row").Value2 = $r.M
    $ws.Range("N$row").Value2 = $r.N
    $ws.Range("O$row").Value2 = "Provincia de Limarí"
    $ws.Range("P$row").Value2 = $r.P
    $ws.Range("Q$row").Value2 = $r.Q
    $ws.Range("R$row").Value2 = "Hortaliza"
}
